# Natmi following Dr Hou advice
# Update Ligand-/Receptor-expressing cell counts (columns E & K) from 1 to 3
# for every data row, and recompute the dependent expression / specificity
# metrics (columns G-J, M-T) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (ECs -> Gnai2/C5ar1 -> ECs)
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 150.0354306666667
$ws.Cells.Item(2, 8).Value = 450.106292
$ws.Cells.Item(2, 9).Value = 0.4152507364956075
$ws.Cells.Item(2, 10).Value = 0.4152507364956075
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.193104333333333
$ws.Cells.Item(2, 14).Value = 3.579313
$ws.Cells.Item(2, 15).Value = 0.03883297235786565
$ws.Cells.Item(2, 16).Value = 0.03883297235786565
$ws.Cells.Item(2, 17).Value = 179.0079224819329
$ws.Cells.Item(2, 18).Value = 1611.071302337396
$ws.Cells.Item(2, 19).Value = 0.01612542037191728
$ws.Cells.Item(2, 20).Value = 0.01612542037191728

# Row 3 (ECs -> Gnai2/C5ar1 -> FAPs)
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 150.0354306666667
$ws.Cells.Item(3, 8).Value = 450.106292
$ws.Cells.Item(3, 9).Value = 0.4152507364956075
$ws.Cells.Item(3, 10).Value = 0.4152507364956075
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.530898
$ws.Cells.Item(3, 14).Value = 88.59269400000001
$ws.Cells.Item(3, 15).Value = 0.9611670276421344
$ws.Cells.Item(3, 16).Value = 0.9611670276421344
$ws.Cells.Item(3, 17).Value = 4430.680999403406
$ws.Cells.Item(3, 18).Value = 39876.12899463065
$ws.Cells.Item(3, 19).Value = 0.3991253161236902
$ws.Cells.Item(3, 20).Value = 0.3991253161236902

# Row 4 (FAPs -> Gnai2/C5ar1 -> ECs)
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 68.382243
$ws.Cells.Item(4, 8).Value = 205.146729
$ws.Cells.Item(4, 9).Value = 0.1892604742946246
$ws.Cells.Item(4, 10).Value = 0.1892604742946246
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.193104333333333
$ws.Cells.Item(4, 14).Value = 3.579313
$ws.Cells.Item(4, 15).Value = 0.03883297235786565
$ws.Cells.Item(4, 16).Value = 0.03883297235786565
$ws.Cells.Item(4, 17).Value = 81.58715044635301
$ws.Cells.Item(4, 18).Value = 734.284354017177
$ws.Cells.Item(4, 19).Value = 0.007349546766719702
$ws.Cells.Item(4, 20).Value = 0.0073495467667197

# Row 5 (FAPs -> Gnai2/C5ar1 -> FAPs)
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 68.382243
$ws.Cells.Item(5, 8).Value = 205.146729
$ws.Cells.Item(5, 9).Value = 0.1892604742946246
$ws.Cells.Item(5, 10).Value = 0.1892604742946246
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 29.530898
$ws.Cells.Item(5, 14).Value = 88.59269400000001
$ws.Cells.Item(5, 15).Value = 0.9611670276421344
$ws.Cells.Item(5, 16).Value = 0.9611670276421344
$ws.Cells.Item(5, 17).Value = 2019.389043044214
$ws.Cells.Item(5, 18).Value = 18174.50138739793
$ws.Cells.Item(5, 19).Value = 0.1819109275279049
$ws.Cells.Item(5, 20).Value = 0.1819109275279049

# Row 6 (M2 -> Gnai2/C5ar1 -> ECs)
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 104.737245
$ws.Cells.Item(6, 8).Value = 314.211735
$ws.Cells.Item(6, 9).Value = 0.2898796499701289
$ws.Cells.Item(6, 10).Value = 0.2898796499701289
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.193104333333333
$ws.Cells.Item(6, 14).Value = 3.579313
$ws.Cells.Item(6, 15).Value = 0.03883297235786565
$ws.Cells.Item(6, 16).Value = 0.03883297235786565
$ws.Cells.Item(6, 17).Value = 124.962460870895
$ws.Cells.Item(6, 18).Value = 1124.662147838055
$ws.Cells.Item(6, 19).Value = 0.01125688843439779
$ws.Cells.Item(6, 20).Value = 0.01125688843439779

# Row 7 (M2 -> Gnai2/C5ar1 -> FAPs)
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 104.737245
$ws.Cells.Item(7, 8).Value = 314.211735
$ws.Cells.Item(7, 9).Value = 0.2898796499701289
$ws.Cells.Item(7, 10).Value = 0.2898796499701289
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 29.530898
$ws.Cells.Item(7, 14).Value = 88.59269400000001
$ws.Cells.Item(7, 15).Value = 0.9611670276421344
$ws.Cells.Item(7, 16).Value = 0.9611670276421344
$ws.Cells.Item(7, 17).Value = 3092.98489889601
$ws.Cells.Item(7, 18).Value = 27836.86409006409
$ws.Cells.Item(7, 19).Value = 0.2786227615357312
$ws.Cells.Item(7, 20).Value = 0.2786227615357311

# Row 8 (sCs -> Gnai2/C5ar1 -> ECs)
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 38.15794
$ws.Cells.Item(8, 8).Value = 114.47382
$ws.Cells.Item(8, 9).Value = 0.105609139239639
$ws.Cells.Item(8, 10).Value = 0.105609139239639
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.193104333333333
$ws.Cells.Item(8, 14).Value = 3.579313
$ws.Cells.Item(8, 15).Value = 0.03883297235786565
$ws.Cells.Item(8, 16).Value = 0.03883297235786565
$ws.Cells.Item(8, 17).Value = 45.52640356507334
$ws.Cells.Item(8, 18).Value = 409.73763208566
$ws.Cells.Item(8, 19).Value = 0.004101116784830886
$ws.Cells.Item(8, 20).Value = 0.004101116784830885

# Row 9 (sCs -> Gnai2/C5ar1 -> FAPs)
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 38.15794
$ws.Cells.Item(9, 8).Value = 114.47382
$ws.Cells.Item(9, 9).Value = 0.105609139239639
$ws.Cells.Item(9, 10).Value = 0.105609139239639
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 29.530898
$ws.Cells.Item(9, 14).Value = 88.59269400000001
$ws.Cells.Item(9, 15).Value = 0.9611670276421344
$ws.Cells.Item(9, 16).Value = 0.9611670276421344
$ws.Cells.Item(9, 17).Value = 1126.83823403012
$ws.Cells.Item(9, 18).Value = 10141.54410627108
$ws.Cells.Item(9, 19).Value = 0.1015080224548081
$ws.Cells.Item(9, 20).Value = 0.1015080224548081

